$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Version bump 2.0.0 -> 2.0.1 reflected as line-number shifts in the stack trace text.
Replace-Text "M2DocEvaluator.caseQuery(M2DocEvaluator.java:543)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:555)"
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1084)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)"
Replace-Text "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1300)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1305)"
Replace-Text "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:278)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:283)"
Replace-Text "M2DocEvaluator.generate(M2DocEvaluator.java:267)" "M2DocEvaluator.generate(M2DocEvaluator.java:272)"
Replace-Text "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:479)"
Replace-Text "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:384)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:388)"
